$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '98.427.27'
$ws.Range('E2').Value = '  -0.40%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.428.95'
$ws.Range('E3').Value = '  +2.38%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '256.26'
$ws.Range('E5').Value = '  -1.43%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '668.14'
$ws.Range('E6').Value = '  +2.43%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.46'
$ws.Range('E7').Value = '  -5.44%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.438'
$ws.Range('E8').Value = '  -6.35%  '

$ws.Range('E9').Value = '  -2.00%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.425.86'
$ws.Range('E11').Value = '  +2.38%  '

$ws.Range('E12').Value = '  +3.17%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.59'
$ws.Range('E13').Value = '  -2.68%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.50'
$ws.Range('E14').Value = '  +16.29%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '98.239.54'
$ws.Range('E15').Value = '  -0.47%  '

$ws.Range('E16').Value = '  -0.51%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.072.00'
$ws.Range('E17').Value = '  +2.10%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.94'
$ws.Range('E18').Value = '  +18.67%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.425.47'
$ws.Range('E19').Value = '  +2.41%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.567'
$ws.Range('E20').Value = '  +30.47%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.80'
$ws.Range('E21').Value = '  +5.47%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.31'
$ws.Range('E22').Value = '  +9.65%  '

$ws.Range('E23').Value = '  -3.68%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '513.89'
$ws.Range('E24').Value = '  -4.53%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000208'
$ws.Range('E25').Value = '  -2.18%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.71'
$ws.Range('E26').Value = '  +7.09%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '101.85'
$ws.Range('E27').Value = '  -0.98%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.98'
$ws.Range('E28').Value = '  +1.83%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.600.92'
$ws.Range('E29').Value = '  +2.07%  '

$ws.Range('E30').Value = '  +0.93%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.68'
$ws.Range('E31').Value = '  +6.21%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.199'
$ws.Range('E32').Value = '  +3.74%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  +0.02%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.52'
$ws.Range('E34').Value = '  +20.50%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  +0.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.577'
$ws.Range('E36').Value = '  +7.43%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '30.20'
$ws.Range('E37').Value = '  +2.80%  '

$ws.Range('E38').Value = '  +15.84%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.03'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '543.80'
$ws.Range('E40').Value = '  +4.41%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.154'
$ws.Range('E41').Value = '  -0.77%  '

$ws.Range('E42').Value = '  +0.00%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.883'
$ws.Range('E43').Value = '  +6.92%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '24.73'
$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('B45').Value = 'MantraDAO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.80'
$ws.Range('E45').Value = '  +0.84%  '

$ws.Range('B46').Value = 'Cosmos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.01'
$ws.Range('E46').Value = '  +14.31%  '

$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.85'
$ws.Range('E47').Value = '  +14.91%  '

$ws.Range('B48').Value = 'ImmutableX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.74'
$ws.Range('E48').Value = '  +16.82%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0430'
$ws.Range('E49').Value = '  +0.71%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.30'
$ws.Range('E50').Value = '  -2.93%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.09'
$ws.Range('E51').Value = '  +9.51%  '
